$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.04575233333333333
$ws.Range("H2").Value = 0.137257
$ws.Range("I2").Value = 0.14147347546269
$ws.Range("J2").Value = 0.14147347546269
$ws.Range("M2").Value = 0.8366046666666667
$ws.Range("N2").Value = 2.509814
$ws.Range("O2").Value = 0.08025679986157715
$ws.Range("P2").Value = 0.08025679986157715
$ws.Range("Q2").Value = 0.03827661557755555
$ws.Range("R2").Value = 0.344489540198
$ws.Range("S2").Value = 0.01135420840593086
$ws.Range("T2").Value = 0.01135420840593086
$ws.Range("G3").Value = 0.04575233333333333
$ws.Range("H3").Value = 0.137257
$ws.Range("I3").Value = 0.14147347546269
$ws.Range("J3").Value = 0.14147347546269
$ws.Range("M3").Value = 7.939250333333333
$ws.Range("O3").Value = 0.7616247559221037
$ws.Range("P3").Value = 0.7616247559221038
$ws.Range("Q3").Value = 0.3632392276674444
$ws.Range("R3").Value = 3.269153049007
$ws.Range("S3").Value = 0.107749701218723
$ws.Range("T3").Value = 0.107749701218723
$ws.Range("G4").Value = 0.04575233333333333
$ws.Range("H4").Value = 0.137257
$ws.Range("I4").Value = 0.14147347546269
$ws.Range("J4").Value = 0.14147347546269
$ws.Range("M4").Value = 1.648242
$ws.Range("N4").Value = 4.944726
$ws.Range("O4").Value = 0.1581184442163192
$ws.Range("P4").Value = 0.1581184442163192
$ws.Range("Q4").Value = 0.075410917398
$ws.Range("R4").Value = 0.678698256582
$ws.Range("S4").Value = 0.02236956583803615
$ws.Range("T4").Value = 0.02236956583803615
$ws.Range("I5").Value = 0.2599907647526892
$ws.Range("J5").Value = 0.2599907647526892
$ws.Range("M5").Value = 0.8366046666666667
$ws.Range("N5").Value = 2.509814
$ws.Range("O5").Value = 0.08025679986157715
$ws.Range("P5").Value = 0.08025679986157715
$ws.Range("Q5").Value = 0.07034227810977779
$ws.Range("R5").Value = 0.633080502988
$ws.Range("S5").Value = 0.02086602677261496
$ws.Range("T5").Value = 0.02086602677261496
$ws.Range("I6").Value = 0.2599907647526892
$ws.Range("J6").Value = 0.2599907647526892
$ws.Range("M6").Value = 7.939250333333333
$ws.Range("O6").Value = 0.7616247559221037
$ws.Range("P6").Value = 0.7616247559221038
$ws.Range("Q6").Value = 0.6675374608602223
$ws.Range("R6").Value = 6.007837147742
$ws.Range("S6").Value = 0.198015402746768
$ws.Range("T6").Value = 0.198015402746768
$ws.Range("I7").Value = 0.2599907647526892
$ws.Range("J7").Value = 0.2599907647526892
$ws.Range("M7").Value = 1.648242
$ws.Range("N7").Value = 4.944726
$ws.Range("O7").Value = 0.1581184442163192
$ws.Range("P7").Value = 0.1581184442163192
$ws.Range("Q7").Value = 0.138585286188
$ws.Range("R7").Value = 1.247267575692
$ws.Range("S7").Value = 0.04110933523330625
$ws.Range("T7").Value = 0.04110933523330625
$ws.Range("G8").Value = 0.1935656666666667
$ws.Range("H8").Value = 0.580697
$ws.Range("I8").Value = 0.5985357597846208
$ws.Range("J8").Value = 0.5985357597846208
$ws.Range("M8").Value = 0.8366046666666667
$ws.Range("N8").Value = 2.509814
$ws.Range("O8").Value = 0.08025679986157715
$ws.Range("P8").Value = 0.08025679986157715
$ws.Range("Q8").Value = 0.1619379400397778
$ws.Range("R8").Value = 1.457441460358
$ws.Range("S8").Value = 0.04803656468303133
$ws.Range("T8").Value = 0.04803656468303133
$ws.Range("G9").Value = 0.1935656666666667
$ws.Range("H9").Value = 0.580697
$ws.Range("I9").Value = 0.5985357597846208
$ws.Range("J9").Value = 0.5985357597846208
$ws.Range("M9").Value = 7.939250333333333
$ws.Range("O9").Value = 0.7616247559221037
$ws.Range("P9").Value = 0.7616247559221038
$ws.Range("Q9").Value = 1.536766283605222
$ws.Range("R9").Value = 13.830896552447
$ws.Range("S9").Value = 0.4558596519566127
$ws.Range("T9").Value = 0.4558596519566128
$ws.Range("G10").Value = 0.1935656666666667
$ws.Range("H10").Value = 0.580697
$ws.Range("I10").Value = 0.5985357597846208
$ws.Range("J10").Value = 0.5985357597846208
$ws.Range("M10").Value = 1.648242
$ws.Range("N10").Value = 4.944726
$ws.Range("O10").Value = 0.1581184442163192
$ws.Range("P10").Value = 0.1581184442163192
$ws.Range("Q10").Value = 0.319043061558
$ws.Range("R10").Value = 2.871387554022
$ws.Range("S10").Value = 0.09463954314497679
$ws.Range("T10").Value = 0.0946395431449768
